# ISM prices-paid monthly series: a new month's reading (Sep, Oct, Nov, Dec 2024)
# has been prepended to the top of the table (newest-first layout), pushing the
# existing history down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 4 new observations above the current row 2 (old row 2
# becomes row 6, etc.)
$ws.Rows("2:5").Insert()

# The inserted rows don't inherit the date number format automatically -
# reapply it so column A keeps displaying as m/d/yy, matching the rest of
# the column.
$ws.Range("A2:A5").NumberFormat = "m/d/yy"

# Newest-first: 12/31/24, 11/30/24, 10/31/24, 9/30/24
$newDates = @(45657, 45626, 45596, 45565)
$newValues = @(64.400000000000006, 58.2, 58.1, 59.4)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newValues[$i]
}

# Column A is a bestFit/autosized column; the new dates ("12/31/24" etc.) are
# a character wider than the widest previous entry, so the column widens
# slightly to keep fitting its content.
$ws.Columns("A").ColumnWidth = 9.6
